# Dungeon.xlsx edit: optimise the count method of enemy and elite
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Background image name for the "losttrees" dungeon entry (row 4) renamed
$ws.Range("P4").Value = "lostftrees"

# QuestDungeon monster/encounter lists (column N) - rename some encounter keys
$ws.Range("N4").Value = "fight;7|fighte;2|ftrees;2|emanflower;2|river;2|cliff;2|losttree;1|oldtree;1|cardbot;2|cardshop;2"
$ws.Range("N5").Value = "fight;5|fighte;2|ftrees;4"
$ws.Range("N6").Value = "fight;10|fighte;3|fsandland;2|potteryroom;2|fhoneyhome;2|esnare;1|basement;1|woodhouse2;1|ebooty;1|trapspear;2|trapdrop;1|potteryman;1|stonedoor2;1|crystalball;2|cardbot;2|cardshop;2"
$ws.Range("N7").Value = "fight;5|fighte;2|ftrees;4"

# Move the active selection to E6
$ws.Range("E6").Select()
